# Fix print functionality for windows with office converter php module
# Applies text substitutions to the wedding-dispensation letter template
# (groom/bride bio fields + ceremony date/time/place + issue date).

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $rng.Find.Execute($find, $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0) | Out-Null
    if ($rng.Find.Found) {
        $rng.Text = $replace
    }
}

# --- CALON SUAMI (groom) block ---
Replace-Text "Squidward" "Ayama Sukuna"

# Longer/more specific string first, so the shorter "Bikini Bottom, " run
# used for the groom's birthplace is not matched by accident.
Replace-Text "Bikini Bottom, 22 Mei 2002" "Japan, 09 Desember 2000"
Replace-Text "Bikini Bottom, " "Japan, "
Replace-Text "08 Juli 2000" "03 September 1998"

Replace-Text "Kasir Krusty krab" "Penyihir"
Replace-Text "Belum Menikah " "Belum Menikah"
Replace-Text "Jl. Pahlawan, Alun-alun Contong, Kec. Bubutan, Surabaya, Jawa Timur 60174" "Jalan Kondang"

# --- CALON ISTRI (bride) block ---
Replace-Text "Sandy Ciks" "Sayumi"
Replace-Text "Buddha" "Islam"
Replace-Text "Pelatih karate" "Test"
Replace-Text "Sudah Menikah" "Belum Menikah"
Replace-Text "Dsn Jaten, Kali Jaten, Selotapak, Trawas, Mojokerto Regency, East Java 61375" "Test"

# --- Ceremony details ---
Replace-Text "senin," "Jum'at,"
Replace-Text "tanggal 11 Februari 2025" "tanggal 02 September 2000"
Replace-Text "11.30 WIB bertempat di" "09.00 WIB bertempat di"
Replace-Text "Rumah kaca," "Shibuya,"

# --- Issue date ---
Replace-Text " 07 November 2024" " 08 November 2024"
